$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value2 = 6509.2
$ws.Range("J70").Value2 = 6509.2
$ws.Range("L70").Value2 = 19527.6
$ws.Range("N70").Value2 = -20067.6
$ws.Range("H73").Value2 = 6509.2
$ws.Range("J73").Value2 = 6509.2
$ws.Range("L73").Value2 = 19527.6
$ws.Range("N73").Value2 = -21399.6
$ws.Range("H92").Value2 = 1364.3334
$ws.Range("I92").Value2 = 1446.25
$ws.Range("J92").Value2 = 1200.5
$ws.Range("K92").Value2 = 1446.25
$ws.Range("L92").Value2 = 1200.5
$ws.Range("M92").Value2 = -198.25
$ws.Range("N92").Value2 = -3696.5
$ws.Range("H99").Value2 = 62507190
$ws.Range("J99").Value2 = 142872740
$ws.Range("L99").Value2 = 428618220
$ws.Range("N99").Value2 = -428621216
$ws.Range("H138").Value2 = 5848.9546
$ws.Range("I138").Value2 = 2459.3076
$ws.Range("J138").Value2 = 7270.4194
$ws.Range("K138").Value2 = 7377.9228
$ws.Range("L138").Value2 = 21811.2582
$ws.Range("M138").Value2 = -2237.9228
$ws.Range("N138").Value2 = -32091.2582

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 5410.794
$ws.Range("I61").Value2 = 4999.0356
$ws.Range("J61").Value2 = 7332.3335
$ws.Range("K61").Value2 = 4999.0356
$ws.Range("L61").Value2 = 7332.3335
$ws.Range("M61").Value2 = -4787.0356
$ws.Range("N61").Value2 = -7756.3335
$ws.Range("H136").Value2 = 5410.794
$ws.Range("I136").Value2 = 4999.0356
$ws.Range("J136").Value2 = 7332.3335
$ws.Range("K136").Value2 = 14997.1068
$ws.Range("L136").Value2 = 21997.0005
$ws.Range("M136").Value2 = -12447.1068
$ws.Range("N136").Value2 = -27097.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value2 = 3911.3333
$ws.Range("I64").Value2 = 871.6
$ws.Range("J64").Value2 = 6082.5713
$ws.Range("K64").Value2 = 871.6
$ws.Range("L64").Value2 = 6082.5713
$ws.Range("M64").Value2 = -646.6
$ws.Range("N64").Value2 = -6532.5713
$ws.Range("H67").Value2 = 3911.3333
$ws.Range("I67").Value2 = 871.6
$ws.Range("J67").Value2 = 6082.5713
$ws.Range("K67").Value2 = 871.6
$ws.Range("L67").Value2 = 6082.5713
$ws.Range("M67").Value2 = -91.60000000000002
$ws.Range("N67").Value2 = -7642.5713
$ws.Range("H134").Value2 = 2176111
$ws.Range("I134").Value2 = 2655037
$ws.Range("J134").Value2 = 20944
$ws.Range("K134").Value2 = 7965111
$ws.Range("L134").Value2 = 62832
$ws.Range("M134").Value2 = -7962576
$ws.Range("N134").Value2 = -67902

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 45460350
$ws.Range("I16").Value2 = 100003950
$ws.Range("J16").Value2 = 7349
$ws.Range("K16").Value2 = 100003950
$ws.Range("L16").Value2 = 7349
$ws.Range("M16").Value2 = -100003663
$ws.Range("N16").Value2 = -7923
$ws.Range("H113").Value2 = 45460350
$ws.Range("I113").Value2 = 100003950
$ws.Range("J113").Value2 = 7349
$ws.Range("K113").Value2 = 100003950
$ws.Range("L113").Value2 = 7349
$ws.Range("M113").Value2 = -100001780
$ws.Range("N113").Value2 = -11689
$ws.Range("H122").Value2 = 1575.8518
$ws.Range("I122").Value2 = 1175.5714
$ws.Range("J122").Value2 = 2976.8333
$ws.Range("K122").Value2 = 3526.7142
$ws.Range("L122").Value2 = 8930.499899999999
$ws.Range("M122").Value2 = -1076.7142
$ws.Range("N122").Value2 = -13830.4999
$ws.Range("H132").Value2 = 4845.4814
$ws.Range("I132").Value2 = 4764.9546
$ws.Range("K132").Value2 = 14294.8638
$ws.Range("M132").Value2 = -11764.8638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value2 = 2942780
$ws.Range("I4").Value2 = 3448898.2
$ws.Range("K4").Value2 = 10346694.6
$ws.Range("M4").Value2 = -10346582.6
$ws.Range("H37").Value2 = 200677.55
$ws.Range("J37").Value2 = 200677.55
$ws.Range("L37").Value2 = 602032.6499999999
$ws.Range("N37").Value2 = -602256.6499999999
$ws.Range("H56").Value2 = 7456
$ws.Range("I56").Value2 = 7456
$ws.Range("K56").Value2 = 7456
$ws.Range("M56").Value2 = -6926
$ws.Range("H116").Value2 = 1259
$ws.Range("J116").Value2 = 1000
$ws.Range("L116").Value2 = 3000
$ws.Range("N116").Value2 = -9884

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value2 = 38464664
$ws.Range("I132").Value2 = 66668230
$ws.Range("J132").Value2 = 5254
$ws.Range("K132").Value2 = 200004690
$ws.Range("L132").Value2 = 15762
$ws.Range("M132").Value2 = -200002160
$ws.Range("N132").Value2 = -20822

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value2 = 7498.8335
$ws.Range("I40").Value2 = 4997.6665
$ws.Range("K40").Value2 = 4997.6665
$ws.Range("M40").Value2 = -4861.6665
$ws.Range("H46").Value2 = 50001424
$ws.Range("J46").Value2 = 55557030
$ws.Range("L46").Value2 = 55557030
$ws.Range("N46").Value2 = -55557406
$ws.Range("H87").Value2 = 34665
$ws.Range("J87").Value2 = 26997.5
$ws.Range("L87").Value2 = 26997.5
$ws.Range("N87").Value2 = -29243.5
$ws.Range("H88").Value2 = 21797.25
$ws.Range("I88").Value2 = 22000
$ws.Range("J88").Value2 = 21594.5
$ws.Range("K88").Value2 = 22000
$ws.Range("L88").Value2 = 21594.5
$ws.Range("M88").Value2 = -21572
$ws.Range("N88").Value2 = -22450.5
$ws.Range("H90").Value2 = 34665
$ws.Range("J90").Value2 = 26997.5
$ws.Range("L90").Value2 = 80992.5
$ws.Range("N90").Value2 = -92224.5
$ws.Range("H91").Value2 = 21797.25
$ws.Range("I91").Value2 = 22000
$ws.Range("J91").Value2 = 21594.5
$ws.Range("K91").Value2 = 22000
$ws.Range("L91").Value2 = 21594.5
$ws.Range("M91").Value2 = -20518
$ws.Range("N91").Value2 = -24558.5
$ws.Range("H132").Value2 = 6545.4
$ws.Range("I132").Value2 = 6279.522
$ws.Range("J132").Value2 = 7419
$ws.Range("K132").Value2 = 18838.566
$ws.Range("L132").Value2 = 22257
$ws.Range("M132").Value2 = -16308.566
$ws.Range("N132").Value2 = -27317

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value2 = 1308.2693
$ws.Range("I107").Value2 = 1491.619
$ws.Range("J107").Value2 = 538.2
$ws.Range("K107").Value2 = 4474.857
$ws.Range("L107").Value2 = 1614.6
$ws.Range("M107").Value2 = -2554.857
$ws.Range("N107").Value2 = -5454.6
$ws.Range("H122").Value2 = 2148
$ws.Range("I122").Value2 = 2133.3333
$ws.Range("K122").Value2 = 6399.999899999999
$ws.Range("M122").Value2 = -3949.999899999999
$ws.Range("H136").Value2 = 8339984.5
$ws.Range("I136").Value2 = 10420962
$ws.Range("J136").Value2 = 16074.833
$ws.Range("K136").Value2 = 31262886
$ws.Range("L136").Value2 = 48224.499
$ws.Range("M136").Value2 = -31260336
$ws.Range("N136").Value2 = -53324.499
